$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Lowercase only the ASCII A-Z letters of a string, leaving every other
# character (accented letters, punctuation, control characters such as
# the form-feed stored as _x000C_, etc.) untouched.
function ConvertTo-AsciiLower {
    param([string]$text)
    $chars = $text.ToCharArray()
    for ($pos = 0; $pos -lt $chars.Length; $pos++) {
        $code = [int]$chars[$pos]
        if ($code -ge 65 -and $code -le 90) {
            $chars[$pos] = [char]($code + 32)
        }
    }
    # Build the result via StringBuilder - repeated "+" string concatenation
    # in this runtime can mis-coerce strings that look numeric (e.g. "inf")
    # into a Double, so avoid it entirely.
    $sb = New-Object System.Text.StringBuilder
    for ($pos = 0; $pos -lt $chars.Length; $pos++) {
        $sb.Append($chars[$pos]) | Out-Null
    }
    return $sb.ToString()
}

$usedRange = $ws.UsedRange
$rowCount = $usedRange.Rows.Count
$colCount = $usedRange.Columns.Count
$rowStart = $usedRange.Row
$colStart = $usedRange.Column

for ($ri = 0; $ri -lt $rowCount; $ri++) {
    for ($ci = 0; $ci -lt $colCount; $ci++) {
        $rowIndex = $rowStart + $ri
        $colIndex = $colStart + $ci
        $cell = $ws.Cells.Item($rowIndex, $colIndex)
        $original = $cell.Value2
        if ($original -ne $null -and $original -ne "") {
            $cell.Value = ConvertTo-AsciiLower $original
        }
    }
}
